# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45961

$ws.Range("B2").Value = 69.26000000000001
$ws.Range("C2").Value = 54.36
$ws.Range("D2").Value = 52.27
$ws.Range("E2").Value = 50.81
$ws.Range("F2").Value = 50.4
$ws.Range("G2").Value = 50.81
$ws.Range("H2").Value = 67.48999999999999
$ws.Range("I2").Value = 77.44
$ws.Range("J2").Value = 75.48999999999999
$ws.Range("K2").Value = 58.9
$ws.Range("L2").Value = 34.09
$ws.Range("M2").Value = 18.12
$ws.Range("N2").Value = 17.07
$ws.Range("O2").Value = 15.92
$ws.Range("P2").Value = 15.7
$ws.Range("Q2").Value = 21.75
$ws.Range("R2").Value = 42.3
$ws.Range("S2").Value = 65.81999999999999
$ws.Range("T2").Value = 83.23
$ws.Range("U2").Value = 90.36
$ws.Range("V2").Value = 93.73999999999999
$ws.Range("W2").Value = 87.8
$ws.Range("X2").Value = 77.42
$ws.Range("Y2").Value = 68.48999999999999
$ws.Range("Z2").Value = 55.79

$ws.Range("AB2").Value = 81.86
$ws.Range("AD2").Value = 90.77
$ws.Range("AF2").Value = 86.8
$ws.Range("AG2").Value = "1h-16h"
